$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> asisstencia)
$ws.Name = "asisstencia"

# New attendance column E for the week of 2025-05-15 (serial 45792),
# mirroring the existing D column's date style
$ws.Range("E1").Value = 45792
$ws.Range("E1").NumberFormat = $ws.Range("D1").NumberFormat

# Attendance marks for the new date - default everyone to Present ("P")
$ws.Range("E2:E23").Value = "P"

# Mark absentees ("A") for this date
$ws.Range("E7").Value = "A"
$ws.Range("E20").Value = "A"
$ws.Range("E21").Value = "A"

# Summary row: count of "P" for the new column, same pattern as column D
$ws.Range("E24").Formula = "=COUNTIF(E2:E23,""P"")"
$ws.Rows.Item(24).RowHeight = 15

# Match the selection left by the editing session
$ws.Range("E20").Select() | Out-Null
